$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlinks attached to C5/D5 (the shmulmaor2 / vikicrestina
# e-mail addresses) before the row itself disappears.
$changed = $true
while ($changed) {
    $changed = $false
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$C$5' -or $addr -eq '$D$5') {
            $h.Delete()
            $changed = $true
            break
        }
    }
}

# Delete the entire 5th row (the shmulmaor2 / vikicrestina review entry).
# This shifts the previously-empty row 6 up to become row 5 and removes
# the now-unreferenced shared strings on save.
$ws.Rows("5:5").Delete()

# Move the active selection to A5, matching the post-edit worksheet state.
$ws.Range("A5").Select()
